$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 809
$ws.Range("B2").Value = 983
$ws.Range("C2").Value = 689
$ws.Range("D2").Value = 644
$ws.Range("E2").Value = 434
$ws.Range("F2").Value = 982
$ws.Range("G2").Value = 933
